$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain numeric cells ---
$ws.Range("A7").Value = 112222974
$ws.Range("B7").Value = 57588
$ws.Range("E7").Value = 208242
$ws.Range("Q7").Value = 636002
$ws.Range("R7").Value = 6520657
$ws.Range("S7").Value = 10

# --- Plain text cells (none of these are numeric/date-like, so Excel's
#     automatic type inference leaves them as text) ---
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "LC"
$ws.Range("F7").Value = "Mindre vattensalamander"
$ws.Range("G7").Value = "Lissotriton vulgaris"
$ws.Range("H7").Value = "(Linnaeus, 1758)"
$ws.Range("J7").Value = "ex."
$ws.Range("K7").Value = "adult"
$ws.Range("L7").Value = "hona"
$ws.Range("M7").Value = "i vatten/simmande"
$ws.Range("N7").Value = "observerad"
$ws.Range("P7").Value = "Nynäs, Srm"
$ws.Range("T7").Value = "Södermanland"
$ws.Range("U7").Value = "Nyköping"
$ws.Range("V7").Value = "Södermanland"
$ws.Range("W7").Value = "Bälinge"
$ws.Range("Z7").Value = "21:30"
$ws.Range("AB7").Value = "21:30"
$ws.Range("AW7").Value = "Stefan Andersson"
$ws.Range("AX7").Value = "Stefan Andersson"

# --- Text cells whose literal content looks numeric/date-like: entering
#     them with a plain .Value assignment would let Excel auto-convert
#     them (e.g. "2023-04-18" -> a date serial, "1" -> the number 1).
#     Force them to land as genuine text by writing a TEXT() formula and
#     then flattening the formula to its literal value via copy/paste. ---
$ws.Range("I7").Formula = "=TEXT(1,""0"")"
$ws.Range("I7").Copy()
$ws.Range("I7").PasteSpecial(-4163)

$ws.Range("Y7").Formula = "=TEXT(DATE(2023,4,18),""yyyy-mm-dd"")"
$ws.Range("Y7").Copy()
$ws.Range("Y7").PasteSpecial(-4163)

$ws.Range("AA7").Formula = "=TEXT(DATE(2023,4,18),""yyyy-mm-dd"")"
$ws.Range("AA7").Copy()
$ws.Range("AA7").PasteSpecial(-4163)

# --- Boolean cells ---
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false

$ws.Application.CutCopyMode = $false
